$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$tbl = $ws.ListObjects.Add(1, $ws.Range("A1:U65"), $null, 1)
Write-Output $tbl.ShowHeaders
Write-Output $tbl.ShowTotals
Write-Output $tbl.ShowAutoFilterDropDown
